$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Ativacao date: "01/01/2012" -> "01/01/2023" (B8/C8)
#    Entering a dd/mm/yyyy-looking string directly gets auto-parsed into
#    a date serial by the "smart" input layer, which would also bump the
#    cell's style (new numFmt). Instead compute the text via TEXT()/DATE()
#    then paste back as a value so the result lands back as plain text in
#    the existing (unmodified) cell style.
$ws.Range("B8:C8").Formula = '=TEXT(DATE(2023,1,1),"dd/mm/yyyy")'
$ws.Range("B8:C8").Copy()
$ws.Range("B8:C8").PasteSpecial(-4163)

# ---------------------------------------------------------------------
# 2. Row 10 (Objetivos:) responsible teacher changes
$ws.Range("B10").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C10").Value = "5840730 - Antonio Jefferson da Silva Machado"

# ---------------------------------------------------------------------
# Helper pattern for newly populated B/C cells (previously empty): copy
# the column's established formatting from a neighbouring row first, so
# the new cell lands on the same style as the rest of the column instead
# of whatever default the engine would otherwise pick for a blank cell.

# 3. Row 11 (Objectives:) gains English objectives text
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("B11").Value = "Provide the student with an overview of the area of Solid State Physics, with emphasis on fundamental ideas and general concepts, such as electron gas, elementary excitations, band structure, etc. The course should be rich in experimental results that illustrate general principles and behaviors of solids (eg, behavior of physical quantities with temperature)."
$ws.Range("C11").Value = "Provide the student with an overview of the area of Solid State Physics, with emphasis on fundamental ideas and general concepts, such as electron gas, elementary excitations, band structure, etc. The course should be rich in experimental results that illustrate general principles and behaviors of solids (eg, behavior of physical quantities with temperature)."

# 4. Row 13 (Programa resumido:) responsible teacher changes
$ws.Range("B13").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C13").Value = "5840730 - Antonio Jefferson da Silva Machado"

# 5. Row 14 (Short syllabus:) gains short syllabus text
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("B14").Value = "Crystal structure and bonds. Lattice vibrations, phonons and thermal properties. Free electron Fermi gas. Power bands. Semiconductors. Fermi metals and surfaces."
$ws.Range("C14").Value = "Crystal structure and bonds. Lattice vibrations, phonons and thermal properties. Free electron Fermi gas. Power bands. Semiconductors. Fermi metals and surfaces."

# 6. Row 15 (Programa:) responsible teacher changes
$ws.Range("B15").Value = "5840726 - Cristina Bormio Nunes"
$ws.Range("C15").Value = "5840726 - Cristina Bormio Nunes"

# 7. Row 16 (Syllabus:) gains syllabus text
$ws.Range("B15").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("B16").Value = "¨ Structure of crystals.¨ Crystal diffraction and the reciprocal lattice.¨ Bonds in crystals: ionic crystals and covalent crystals¨ Elastic constants and elastic waves.¨ Crystal vibrations. phonons¨ Fermi gas: free electron model; movement in magnetic fields.¨ Energy bands. Bloch functions.¨ Semiconductor crystals."
$ws.Range("C16").Value = "¨ Structure of crystals.¨ Crystal diffraction and the reciprocal lattice.¨ Bonds in crystals: ionic crystals and covalent crystals¨ Elastic constants and elastic waves.¨ Crystal vibrations. phonons¨ Fermi gas: free electron model; movement in magnetic fields.¨ Energy bands. Bloch functions.¨ Semiconductor crystals."

# ---------------------------------------------------------------------
# 8. Row 20 (Norma de recuperacao value) wording change
$ws.Range("B20").Value = "Média aritmética de duas provas com mesmo peso."
$ws.Range("C20").Value = "Média aritmética de duas provas com mesmo peso."
